# PantryPal Testing.xlsx - add unit testing results to the "Testing" sheet.
# Rows 7-22 (TestID 5-20) get filled in with the remaining test cases that
# were completed: BackEnd tests, Pantry tests, Recipe tests and Shopping
# Cart tests, all marked Pass / Complete with the relevant testers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")

# Columns: B=TestID, C=TestDescription, D=Test Status, E=Testers Assigned,
#          F=Testers Who Completed, G=RequirementID, H=Requirement Status,
#          I=Developers Assigned, J=Developers Who Completed
$rows = @(
    @(5,  "BackEnd - PostgreSQL",                              "Pass", "Elijah, Makaela", "Elijah",  5,  "Complete", "Denver",  "Denver"),
    @(6,  "BackEnd - Spring Boot",                              "Pass", "Elijah",          "Elijah",  6,  "Complete", "Denver",  "Denver"),
    @(7,  "Recipe - Invalid form submission checks",            "Pass", "Elijah",          "Elijah",  7,  "Complete", "Makaela", "Denver"),
    @(8,  "Recipe - Duplicate Recipe check",                    "Pass", "Elijah",          "Elijah",  8,  "Complete", "Makaela", "Denver"),
    @(9,  "Recipe - Delete Recipe",                             "Pass", "Elijah ",         "Denver",  9,  "Complete", "Makaela", "Denver"),
    @(10, "Recipe - Auto remove ingredients from pantry",       "Pass", "Elijah",          "Elijah",  10, "Complete", "Makaela", "Denver"),
    @(11, "Pantry - Invalid form submission checks",            "Pass", "Elijah",          "Elijah",  11, "Complete", "Elijah ", "Denver"),
    @(12, "Pantry - Duplicate ingredients check",                "Pass", "Elijah",          "Elijah",  12, "Complete", "Elijah ", "Denver"),
    @(13, "Pantry - Request for invalid ingredient from Recipe", "Pass", "Elijah",          "Elijah",  13, "Complete", "Elijah ", "Denver"),
    @(14, "Pantry - Delete Ingredient",                          "Pass", "Elijah ",         "Denver",  14, "Complete", "Elijah ", "Denver"),
    @(15, "Pantry - Updates with shopping cart changes",         "Pass", "Elijah",          "Elijah",  15, "Complete", "Elijah",  "Denver"),
    @(16, "Shopping Cart - Invalid manual entry check",          "Pass", "Elijah",          "Elijah",  16, "Complete", "Denver",  "Denver"),
    @(17, "Shopping Cart - Manual Entry",                        "Pass", "Elijah",          "Elijah",  17, "Complete", "Denver",  "Denver"),
    @(18, "Shopping Cart - Auto loading ingredients from Recipe","Pass", "Elijah",          "Elijah",  18, "Complete", "Denver",  "Denver"),
    @(19, "Shopping Cart - Add ingredients to pantry when checkout","Pass","Elijah",        "Elijah",  19, "Complete", "Denver",  "Denver"),
    @(20, "Shopping Cart - Delete Ingredient",                   "Pass", "Elijah",          "Denver",  20, "Complete", "Denver",  "Denver")
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
    $r++
}

# Select I11:J11 to match the saved selection state.
$ws.Range("I11:J11").Select()
